# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the dfcc41c5-... row on the zh-cn and de-de sheets, and
# roll the "Latest HO Xliff Generate Date" on the Overview sheet forward to
# match the newest of the two (de-de's, which is later).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# zh-cn: row 3 is the dfcc41c5-... file.
$zhcn.Range("H3").Value = "2016-08-29 14:55:56"
$zhcn.Range("K3").Value = "2016-08-29 14:56:30"

# de-de: row 3 is the dfcc41c5-... file.
$dede.Range("H3").Value = "2016-08-29 14:56:03"
$dede.Range("K3").Value = "2016-08-29 14:56:38"

# Overview: row 3 is the dfcc41c5-... file; its "Latest HO Xliff Generate
# Date" becomes the newest handoff datetime across languages (de-de's).
$overview.Range("G3").Value = "2016-08-29 14:56:03"
